{"js": "// Update the date heading and all the \"two-digit \u00f7 one-digit\" answer\n// cells in the practice-sheet table to the new day's values.\n//\n// The document body starts with a centered date paragraph, followed by\n// a single 20-row x 5-column table. Only every 4th row (0, 4, 8, 12, 16)\n// actually holds answers; the rows between are blank spacer rows.\n\n// --- 1. Update the date heading paragraph -------------------------------\nconst headingHits = context.document.body.search(\"2025-03-05 Wednesday\", {\n  matchCase: true,\n});\nheadingHits.load(\"items\");\nawait context.sync();\n\nif (headingHits.items.length > 0) {\n  headingHits.items[0].insertText(\"2025-03-06 Thursday\", \"Replace\");\n}\nawait context.sync();\n\n// --- 2. Update the answer cells in the table -----------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col) -> new answer text (rows 0, 4, 8, 12, 16 hold the five\n// \"divide this row\" answers; other rows are blank spacers left untouched).\nconst cellUpdates = [\n  // row 0\n  { row: 0, col: 0, newText: \"77\u00f72=38, 1\" },\n  { row: 0, col: 1, newText: \"57\u00f76=9, 3\" },\n  { row: 0, col: 2, newText: \"77\u00f77=11, 0\" },\n  { row: 0, col: 3, newText: \"82\u00f72=41, 0\" },\n  { row: 0, col: 4, newText: \"28\u00f74=7, 0\" },\n  // row 4\n  { row: 4, col: 0, newText: \"11\u00f76=1, 5\" },\n  { row: 4, col: 1, newText: \"13\u00f75=2, 3\" },\n  { row: 4, col: 2, newText: \"34\u00f76=5, 4\" },\n  { row: 4, col: 3, newText: \"82\u00f78=10, 2\" },\n  { row: 4, col: 4, newText: \"83\u00f74=20, 3\" },\n  // row 8\n  { row: 8, col: 0, newText: \"93\u00f78=11, 5\" },\n  { row: 8, col: 1, newText: \"36\u00f76=6, 0\" },\n  { row: 8, col: 2, newText: \"37\u00f79=4, 1\" },\n  { row: 8, col: 3, newText: \"91\u00f76=15, 1\" },\n  { row: 8, col: 4, newText: \"93\u00f78=11, 5\" },\n  // row 12\n  { row: 12, col: 0, newText: \"39\u00f74=9, 3\" },\n  { row: 12, col: 1, newText: \"57\u00f78=7, 1\" },\n  { row: 12, col: 2, newText: \"77\u00f78=9, 5\" },\n  { row: 12, col: 3, newText: \"14\u00f78=1, 6\" },\n  { row: 12, col: 4, newText: \"41\u00f75=8, 1\" },\n  // row 16\n  { row: 16, col: 0, newText: \"36\u00f75=7, 1\" },\n  { row: 16, col: 1, newText: \"99\u00f73=33, 0\" },\n  { row: 16, col: 2, newText: \"86\u00f74=21, 2\" },\n  { row: 16, col: 3, newText: \"77\u00f78=9, 5\" },\n  { row: 16, col: 4, newText: \"75\u00f75=15, 0\" },\n];\n\nfor (const u of cellUpdates) {\n  table.getCell(u.row, u.col).value = u.newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and all the \"two-digit / one-digit\" answer\n# cells in the practice-sheet table to the new day's values.\n#\n# The document body starts with a centered date paragraph, followed by a\n# single 20-row x 5-column table. Only every 4th row (1, 5, 9, 13, 17 in\n# Word's 1-based Cell() numbering) actually holds answers; the rows in\n# between are blank spacer rows and are left untouched.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date heading paragraph --------------------------------\n$find = $d.Content.Find\n$find.Text = \"2025-03-05 Wednesday\"\n$find.Replacement.Text = \"2025-03-06 Thursday\"\n$find.Execute([ref]\"2025-03-05 Wednesday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-03-06 Thursday\", 2) | Out-Null\n\n# --- 2. Update the answer cells in the table ------------------------------\n$table = $d.Tables.Item(1)\n\n# Word table rows/columns are 1-based: data rows are 1, 5, 9, 13, 17.\n$table.Cell(1, 1).Range.Text = \"77\u00f72=38, 1\"\n$table.Cell(1, 2).Range.Text = \"57\u00f76=9, 3\"\n$table.Cell(1, 3).Range.Text = \"77\u00f77=11, 0\"\n$table.Cell(1, 4).Range.Text = \"82\u00f72=41, 0\"\n$table.Cell(1, 5).Range.Text = \"28\u00f74=7, 0\"\n\n$table.Cell(5, 1).Range.Text = \"11\u00f76=1, 5\"\n$table.Cell(5, 2).Range.Text = \"13\u00f75=2, 3\"\n$table.Cell(5, 3).Range.Text = \"34\u00f76=5, 4\"\n$table.Cell(5, 4).Range.Text = \"82\u00f78=10, 2\"\n$table.Cell(5, 5).Range.Text = \"83\u00f74=20, 3\"\n\n$table.Cell(9, 1).Range.Text = \"93\u00f78=11, 5\"\n$table.Cell(9, 2).Range.Text = \"36\u00f76=6, 0\"\n$table.Cell(9, 3).Range.Text = \"37\u00f79=4, 1\"\n$table.Cell(9, 4).Range.Text = \"91\u00f76=15, 1\"\n$table.Cell(9, 5).Range.Text = \"93\u00f78=11, 5\"\n\n$table.Cell(13, 1).Range.Text = \"39\u00f74=9, 3\"\n$table.Cell(13, 2).Range.Text = \"57\u00f78=7, 1\"\n$table.Cell(13, 3).Range.Text = \"77\u00f78=9, 5\"\n$table.Cell(13, 4).Range.Text = \"14\u00f78=1, 6\"\n$table.Cell(13, 5).Range.Text = \"41\u00f75=8, 1\"\n\n$table.Cell(17, 1).Range.Text = \"36\u00f75=7, 1\"\n$table.Cell(17, 2).Range.Text = \"99\u00f73=33, 0\"\n$table.Cell(17, 3).Range.Text = \"86\u00f74=21, 2\"\n$table.Cell(17, 4).Range.Text = \"77\u00f78=9, 5\"\n$table.Cell(17, 5).Range.Text = \"75\u00f75=15, 0\"\n"}
